$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9, shifting existing rows 9-10 down to 10-11
$ws.Rows.Item(9).Insert()

# Fill in the new row 9 with the new data record
$ws.Cells.Item(9, 1).Value = 11
$ws.Cells.Item(9, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value = "Bíobío"
$ws.Cells.Item(9, 4).Value = 44875
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100103
$ws.Cells.Item(9, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(9, 9).Value = 100103003
$ws.Cells.Item(9, 10).Value = "Damasco"
$ws.Cells.Item(9, 11).Value = "Castle Brite"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 50
$ws.Cells.Item(9, 14).Value = 31000
$ws.Cells.Item(9, 15).Value = 32000
$ws.Cells.Item(9, 16).Value = 31400
$ws.Cells.Item(9, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 19).Value = 3140
$ws.Cells.Item(9, 20).Value = 10
